$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C7").Value = 7007
$ws.Range("E7").Value = 289917074
$ws.Range("C14").Value = 110814
$ws.Range("E14").Value = 253241833
$ws.Range("C37").Value = 23047
$ws.Range("E37").Value = 130182720
$ws.Range("C53").Value = 141676
$ws.Range("E53").Value = 590050995
$ws.Range("C56").Value = 11974
$ws.Range("E56").Value = 187807204
$ws.Range("C63").Value = 14338
$ws.Range("E63").Value = 36137488
$ws.Range("C65").Value = 2010
$ws.Range("E65").Value = 13580554
$ws.Range("C70").Value = 15717
$ws.Range("E70").Value = 24651442
$ws.Range("C74").Value = 938
$ws.Range("E74").Value = 4168041
$ws.Range("C79").Value = 116587
$ws.Range("E79").Value = 447322464
$ws.Range("C81").Value = 17431
$ws.Range("E81").Value = 133553933
$ws.Range("C90").Value = 34348
$ws.Range("E90").Value = 67213232
$ws.Range("C91").Value = 151073
$ws.Range("E91").Value = 481901596
$ws.Range("C92").Value = 408954
$ws.Range("D92").Value = 70903
$ws.Range("E92").Value = 1593050860
$ws.Range("C93").Value = 209443
$ws.Range("E93").Value = 1306950791
$ws.Range("C94").Value = 94125
$ws.Range("E94").Value = 914927107
$ws.Range("C95").Value = 50689
$ws.Range("E95").Value = 928906921
$ws.Range("E96").Value = 786973618
$ws.Range("C97").Value = 2150
$ws.Range("E97").Value = 213846316
$ws.Range("C98").Value = 807
$ws.Range("E98").Value = 117420097
$ws.Range("C104").Value = 135210
$ws.Range("D104").Value = 23286
$ws.Range("E104").Value = 272060324
$ws.Range("C106").Value = 18331
$ws.Range("E106").Value = 41271449
$ws.Range("C107").Value = 6384
$ws.Range("E107").Value = 21931803
$ws.Range("C108").Value = 2827
$ws.Range("E108").Value = 18460239
$ws.Range("C113").Value = 8800
$ws.Range("E113").Value = 12651124
$ws.Range("E114").Value = 9073557
$ws.Range("C115").Value = 11680
$ws.Range("E115").Value = 32886288
$ws.Range("C116").Value = 4548
$ws.Range("E116").Value = 20411989
$ws.Range("C122").Value = 8484
$ws.Range("E122").Value = 12669069
$ws.Range("C131").Value = 75580
$ws.Range("E131").Value = 307202531
$ws.Range("C138").Value = 15
$ws.Range("E138").Value = 626897
$ws.Range("C142").Value = 168966
$ws.Range("E142").Value = 681733793
$ws.Range("C165").Value = 83801
$ws.Range("D165").Value = 17112
$ws.Range("E165").Value = 354961515
$ws.Range("C167").Value = 12217
$ws.Range("E167").Value = 105725619
$ws.Range("C168").Value = 6204
$ws.Range("E168").Value = 100524964
$ws.Range("C174").Value = 226077
$ws.Range("E174").Value = 900500557
$ws.Range("C177").Value = 14706
$ws.Range("E177").Value = 251083376
